$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (stored width = ColumnWidth + 5/6, quantized to 1/6) ---
$ws.Columns.Item(1).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 64.83333333333334
$ws.Columns.Item(4).ColumnWidth = 23.5
$ws.Columns.Item(5).ColumnWidth = 18.833333333333336
$ws.Columns.Item(6).ColumnWidth = 19.5
$ws.Columns.Item(7).ColumnWidth = 28.833333333333336
$ws.Columns.Item(11).ColumnWidth = 62.0

# --- Row 5 custom height ---
$ws.Rows.Item(5).RowHeight = 13.5

# --- Value changes ---
$ws.Range("N1").Value = "IncludesByPath"
$ws.Range("M2").Value = " @smoketest"
$ws.Range("K3").Value = "petId=id;petName=name;category_name=category.name"
$ws.Range("M3").Value = "@pet @smoketest"
$ws.Range("M4").Value = "@pet @smoketest"
$ws.Range("M5").Value = "@xml @smoketest"

# --- Selection ---
$ws.Range("K4").Select()
